$d = $word.ActiveDocument

# --- Paragraph 2: "World Environment: (using raylib library)" -------------
# Bump heading font size to 12pt (sz/szCs = 24 half-points) and split the
# run so "raylib" is wrapped in spell-check proofErr markers, matching the
# target OOXML exactly (paragraph mark rPr + 3 runs).
$p2 = $d.Paragraphs(2)
$xmlWorldEnv = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">World Environment: (using </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t>raylib</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve"> library)</w:t></w:r>' + `
  '</w:p>'
$p2.Range.InsertXML($xmlWorldEnv)

# --- Paragraph 4: "Functionality of the AI:" -------------------------------
$p4 = $d.Paragraphs(4)
$p4.Range.Font.Size = 12
$p4.Range.Font.SizeBi = 12

# --- Paragraph 18: "Interaction with the Simulated Environment" ------------
$p18 = $d.Paragraphs(18)
$p18.Range.Font.Size = 12
$p18.Range.Font.SizeBi = 12

# --- Paragraph 21: "Difficulty Levels and Their Controls" ------------------
$p21 = $d.Paragraphs(21)
$p21.Range.Font.Size = 12
$p21.Range.Font.SizeBi = 12
